$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.923.95'
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").Value = '2.416.42'
$ws.Range("E3").Value = '  +0.33%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '562.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.48%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.83'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.49%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +0.29%  '
$ws.Range("E9").Value = '  +0.58%  '
$ws.Range("E10").Value = '  -1.36%  '
$ws.Range("E11").Value = '  -3.62%  '
$ws.Range("E12").Value = '  -0.37%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '25.66'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.04%  '
$ws.Range("E14").Value = '  -0.46%  '
$ws.Range("D15").Value = '2.853.50'
$ws.Range("E15").Value = '  +0.49%  '
$ws.Range("D16").Value = '61.934.12'
$ws.Range("E16").Value = '  +0.17%  '
$ws.Range("D17").Value = '2.427.47'
$ws.Range("E17").Value = '  +1.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.28'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.87%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '322.65'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("E20").Value = '  +1.97%  '
$ws.Range("E21").Value = '  -1.06%  '
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.50'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.31%  '
$ws.Range("E24").Value = '  +1.66%  '
$ws.Range("E25").Value = '  -3.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '557.42'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.01%  '
$ws.Range("D27").Value = '2.536.43'
$ws.Range("E27").Value = '  +0.43%  '
$ws.Range("E28").Value = '  -0.10%  '
$ws.Range("D29").Value = '0.0₃0932'
$ws.Range("E29").Value = '  +1.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.16'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.85%  '
$ws.Range("E31").Value = '  -3.79%  '
$ws.Range("E32").Value = '  -0.50%  '
$ws.Range("E33").Value = '  -0.85%  '
$ws.Range("E34").Value = '  -3.72%  '
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("E36").Value = '  +0.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.378'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.74%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '153.75'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.89%  '
$ws.Range("E39").Value = '  -3.69%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.50'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.56%  '
$ws.Range("E41").Value = '  -1.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.993'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.71%  '
$ws.Range("E43").Value = '  -2.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '147.19'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0526'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.95%  '
$ws.Range("E47").Value = '  +0.76%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '19.75'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.53%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0919'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.35%  '
$ws.Range("E50").Value = '  -0.38%  '
$ws.Range("E51").Value = '  +0.76%  '
